$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 1.33
$ws.Range("P2").Value = 3.25
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.7
